# Applies:
#  - Metadata!B8 (Date) updated to the new generation timestamp
#  - Elements sheet: remove the "fr-lm-naissance.resultat" row (row 7),
#    shifting subsequent rows (identificationNouveauNe, observationNaissance) up.

$wb = $excel.ActiveWorkbook

# --- Update the Date property on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# --- Remove the "resultat" element row from the Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Rows.Item(7).Delete()
